# Remove the "Device:" / "${device.deviceName}" row contents (A4:B4).
# Clearing the contents (but not the formatting) of those two cells is
# exactly what the xlsx diff shows: A4/B4 keep their cell styles (s="1"
# / s="6") but lose their t="s"/<v> shared-string references, the two
# now-unused shared strings ("Device:" and "${device.deviceName}") drop
# out of sharedStrings.xml, and every subsequent shared-string index
# shifts down by 2. Re-selecting the cleared range also updates the
# sheet's active selection/view to point at A4:B4 instead of the old
# N9 selection (matches the "Remove device Name ... and change sheet
# Name" commit message).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4:B4").ClearContents()

$ws.Range("A4:B4").Select()
